$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Create tables in database and generate scripts" task (row 5) as Done,
# matching the "Done" formatting already used by rows 3 and 4 (E3/E4): copy the
# cell format from E4 onto E5, then set its value to "Done".
$ws.Range("E4").Copy()
$ws.Range("E5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("E5").Value = "Done"

# Update the active cell selection to D10 (matches workbook UI state change in the diff).
$ws.Range("D10").Select()
